$p = $ppt.ActivePresentation
$lf = [char]10

# 1) Delete the 'Glucose production' slide (old slide 4; CO2/Water content now ends at slide 3,
#    and the former slide 5 'Add improvements!' becomes the new slide 4).
$p.Slides.Item(4).Delete()

# 2) Rewrite the speaker notes body for slide 1 (Photosynthesis).
$p.Slides.Item(1).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = '.' + $lf + '' + $lf + 'Plant life, as we know it, thrives on a process known as photosynthesis. This biochemical reaction, which occurs within the plant cells, is responsible for the production of food and oxygen that sustains all life forms on Earth. The process can be simply described as the conversion of carbon dioxide, water, and sunlight into glucose (a type of sugar) and oxygen.' + $lf + '' + $lf + 'The photosynthetic process is divided into two main stages: the light-dependent reactions and the light-independent reactions, also known as Calvin cycle. In the light-dependent reactions, energy from sunlight is captured by pigments such as chlorophyll, found in the thylakoid membranes of chloroplasts. This energy is used to convert water molecules into protons, electrons, and oxygen. The protons create a gradient across the thylakoid membrane, driving the flow of electrons through a series of proteins in a pathway called the electron transport chain. This flow generates ATP (adenosine triphosphate), the energy currency of the cell.' + $lf + '' + $lf + 'The light-independent reactions occur in the stroma of the chloroplasts. Here, the ATP and NADPH (nicotinamide adenine dinucleotide phosphate) produced during the light-dependent reactions provide the energy required to convert carbon dioxide into glucose. This process involves several enzyme-catalyzed reactions, collectively known as the Calvin cycle. During the Calvin cycle, carbon dioxide is fixed into an organic compound, forming a three-carbon molecule called glyceraldehyde 3-phosphate (G3P). G3P can then be converted into glucose or used immediately for energy needs.' + $lf + '' + $lf + 'Photosynthesis is not only essential for plants but also plays a crucial role in maintaining the Earth''s atmosphere. Through photosynthesis, plants absorb carbon dioxide, helping to regulate its levels and prevent excessive accumulation in the atmosphere. Excessive carbon dioxide can lead to increased greenhouse gas emissions, contributing to global warming and climate change. On the other hand, photosynthesis produces oxygen, which is vital for aerobic organisms, including humans, animals, and most microorganisms.' + $lf + '' + $lf + 'In addition to its primary functions, photosynthesis also contributes to the chemical cycling of nutrients in ecosystems. For example, during decomposition, organic matter breaks down, releasing carbon dioxide back into the environment. Plants then absorb this carbon dioxide and convert it back into organic compounds during photosynthesis, effectively recycling carbon back into the ecosystem.' + $lf + '' + $lf + 'However, photosynthesis is not without challenges. Environmental factors such as temperature, light intensity, and carbon dioxide concentration can significantly impact the efficiency of photosynthesis. For instance, high temperatures can inhibit the enzymes involved in the Calvin cycle, reducing the rate of photosynthesis. Similarly, low light intensity or carbon dioxide concentrations can limit the production of ATP and NADPH, respectively, slowing down the photosynthetic process.' + $lf + '' + $lf + 'Researchers are continuously'

# 3) Rewrite the speaker notes body for slide 2 (Carbon dioxide).
$p.Slides.Item(2).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = '' + $lf + '' + $lf + 'Carbon Dioxide (CO2) plays a crucial role in the process of photosynthesis, a fundamental biological activity carried out by green plants, algae, and certain bacteria. This essential chemical compound serves as the primary source of carbon for these organisms, helping them to produce their own food through a series of complex biochemical reactions.' + $lf + '' + $lf + 'In photosynthesis, CO2 is absorbed from the air or water, depending on the organism. Inside the plant cells, within structures called chloroplasts, the CO2 molecules undergo a series of transformations. The initial step involves the carbon dioxide being converted into an organic compound, primarily glucose or other sugars, with the help of sunlight, water, and other nutrients.' + $lf + '' + $lf + 'The process begins when the carbon dioxide molecule interacts with a molecule of ribulose bisphosphate (RuBP), another essential component of photosynthesis. This interaction results in the formation of two molecules of 3-phosphoglycerate, which is then further processed to form glucose or other carbohydrates.' + $lf + '' + $lf + 'This entire process is facilitated by a group of proteins known as Rubisco (Ribulose-1,5-bisphosphate carboxylase/oxygenase). Rubisco is responsible for the carboxylation of RuBP, which is the first step in the conversion of carbon dioxide into organic compounds. However, Rubisco can also catalyze the oxygenation of RuBP, leading to the production of harmful compounds such as glycolaldehyde and phosphoglycolic acid. This side reaction, known as photorespiration, consumes energy and reduces the overall efficiency of photosynthesis.' + $lf + '' + $lf + 'It is worth noting that CO2 is not just a raw material in photosynthesis but also a byproduct of cellular respiration, the opposite process where organisms convert stored carbohydrates back into energy. During cellular respiration, glucose is broken down, releasing CO2 and water, while also generating ATP (adenosine triphosphate), the primary energy currency of cells.' + $lf + '' + $lf + 'The balance between photosynthesis and cellular respiration is critical for the survival of most organisms. Photosynthesis provides the necessary glucose for growth and reproduction, while cellular respiration supplies energy for various cellular functions.' + $lf + '' + $lf + 'The importance of CO2 in photosynthesis extends beyond its role as a source of carbon. CO2 also influences the rate of photosynthesis. An increase in CO2 concentration can stimulate photosynthesis due to the enhanced availability of this essential raw material. This phenomenon, known as the CO2 fertilization effect, has significant implications for global climate change and agriculture.' + $lf + '' + $lf + 'In conclusion, Carbon Dioxide is a key player in the process of photosynthesis, serving as the primary source of carbon for plants, algae, and certain bacteria. Its absorption, conversion into organic compounds, and subsequent release during cellular respiration form a cyclical process that is vital for the survival and growth of these organisms. Understanding the role of CO2 in photosynthesis can provide'

# 4) Rewrite the speaker notes body for slide 3 (Water).
$p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = '' + $lf + '' + $lf + 'The essential role of water in photosynthesis cannot be overstated. This life-sustaining process, carried out by green plants and other organisms, converts carbon dioxide, water, and light energy into glucose and oxygen. The intricate dance of water molecules within this biochemical spectacle is a testament to nature''s ingenuity and efficiency.' + $lf + '' + $lf + 'Water, in its liquid form, serves as the medium for the initial steps of photosynthesis. It is here that water molecules (H2O) are split into hydrogen ions (H+) and electrons (e-) at the thylakoid membrane, during a process known as photophosphorylation. This splitting occurs as a result of the absorption of light energy by chlorophyll, a pigment found within the chloroplasts of plant cells.' + $lf + '' + $lf + 'The energy absorbed by the chlorophyll causes a change in the electronic structure of the chlorophyll molecule, which in turn excites an electron from the water molecule. This excited electron moves through a series of electron transport chains, driving the production of ATP (adenosine triphosphate), the primary energy carrier in cells. Simultaneously, hydrogen ions are pumped across the thylakoid membrane, creating a concentration gradient.' + $lf + '' + $lf + 'This gradient drives the movement of hydrogen ions back into the thylakoid space, a process known as chemiosmosis. As these ions return, they pass through a protein complex called ATP synthase, where their kinetic energy is converted into chemical energy stored in ATP molecules.' + $lf + '' + $lf + 'Meanwhile, in the stroma of the chloroplast, the Calvin cycle takes place. Here, carbon dioxide enters the plant cell through tiny pores called stomata. In the Calvin cycle, carbon dioxide is fixed into organic compounds using the energy stored in ATP and NADPH (nicotinamide adenine dinucleotide phosphate), another reduced coenzyme generated during photophosphorylation.' + $lf + '' + $lf + 'Initially, a five-carbon sugar called ribulose bisphosphate (RuBP) reacts with carbon dioxide to produce two three-carbon sugars, glyceraldehyde 3-phosphate (G3P). These G3P molecules can then be used immediately for energy or stored for later use, such as for growth and reproduction. Additionally, some of the G3P may be used to regenerate RuBP, allowing the Calvin cycle to continue.' + $lf + '' + $lf + 'In essence, water plays a dual role in photosynthesis: it provides the starting material for the process, supplying the carbon needed to build carbohydrates, and it serves as a crucial intermediary, generating the energy required to drive the process forward. The conversion of water into oxygen is a byproduct of this process, making photosynthesis not only essential for life on Earth but also responsible for maintaining our planet''s atmosphere.' + $lf + '' + $lf + 'Understanding the role of water in photosynthesis sheds light on the intricate balance between the environment and living organisms. It underscores the importance of water'

# 5) Update the on-slide title text boxes to match the new, shorter titles.
$p.Slides.Item(2).Shapes.Item(3).TextFrame.TextRange.Text = 'Carbon dioxide'
$p.Slides.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = 'Water'

